$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "authors_short" column for the publications that were added
# from PubMed (rows 9-22), using the "Lastname et al." convention already used above.
$ws.Range("C9").Value  = "Pham et al."
$ws.Range("C10").Value = "Gottwein et al."
$ws.Range("C11").Value = "Hezode et al."
$ws.Range("C12").Value = "Poordad et al."
$ws.Range("C13").Value = "Ng et al."
$ws.Range("C14").Value = "Bourlière et al."
$ws.Range("C15").Value = "Gane et al."
$ws.Range("C16").Value = "Kwo et al."
$ws.Range("C17").Value = "Ng et al."
$ws.Range("C18").Value = "Poordad et al."
$ws.Range("C19").Value = "Lawitz et al."
$ws.Range("C20").Value = "Poordad et al."
$ws.Range("C21").Value = "Foster et al."
$ws.Range("C22").Value = "Curry et al."

# Rename header: "author_short" -> "authors_short" (column C header already existed with the
# old name; just overwrite the text).
$ws.Range("C1").Value = "authors_short"

# Leave the selection where the author ended up after editing the header cell.
$ws.Range("C1").Select()
